# Changed SCMP scraper to work with new BaseScraper changes
# -> "target" is no longer part of the companies list, so its row is
#    removed from the sheet (and the now-unused shared string disappears
#    from sharedStrings.xml on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "target" lives in A21 ("Companies list" header is row 1, so row 21 is
# the 20th company). Deleting the entire row shifts unicredit/vertex/
# walmart up one row, shrinking the sheet from A1:C24 to A1:C23.
$ws.Rows.Item(21).Delete()

# Leave the selection where the author's cursor ended up after the
# deletion - now sitting on "unicredit", which took over A21.
$ws.Range("A21").Select()
